# Move the table legend/description into a new merged title row above the
# header, and freeze panes so the header row stays visible when scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at the top (this shifts the existing header + all data
#    rows down by one, which is exactly what the target diff shows: every
#    row in the sheet moves from r=N to r=N+1).
$ws.Rows.Item(1).Insert()

# 2. Merge A1:J1 to host the new legend text, and give it some breathing
#    room (row height ~39pt, matching the two lines of wrapped legend text).
[void]$ws.Range("A1:J1").Merge()
$ws.Rows.Item(1).RowHeight = 39

# 3. Write the legend text: a bold title run followed by a normal body run
#    describing the methodology (this becomes a rich-text shared string).
$title = "Element gene ontology enrichment"
$body = "`nHPO genes for each element were tested for enrichment among genes co-annotated for gene ontology (GO) terms (hypergeometric test). Bonferroni correction is included as a column, treating each GO term as an independent test."
$full = $title + $body

$ws.Range("A1").Value = $full
$ws.Range("A1").Characters(1, $title.Length).Font.Bold = $true
$ws.Range("A1").Characters($title.Length + 1, $body.Length).Font.Bold = $false

# 4. Wrap the legend text within the merged title cell.
$ws.Range("A1").WrapText = $true

# 5. Freeze panes below the new title row so the column header row (now
#    row 2) stays pinned while scrolling through the data.
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Leave the row-1 range selected across the full row width, matching the
# "select the frozen title row" state captured in the saved view.
[void]$ws.Range("A1:XFD1").Select()
